$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values in row 75 (quarter 01-01-2021) ---
$row75Updates = @{
    "C"  = 9966
    "F"  = 919
    "H"  = 202
    "K"  = 58995
    "M"  = 2121
    "N"  = 7289
    "P"  = 1449
    "R"  = 6552
    "T"  = 9238
    "U"  = 10510
    "W"  = 1871
    "Y"  = 12
    "AD" = 2861
    "AN" = 70514
    "AP" = 10009
    "AU" = 13692
    "AW" = 32815
    "AX" = 55560
    "BA" = 3543
    "BC" = 5127
    "BD" = 689
    "BF" = 4242
}

foreach ($col in $row75Updates.Keys) {
    $ws.Range("$col`75").Value = $row75Updates[$col]
}

# --- Append a new row 76 for quarter 01-04-2021 ---
# Force the value to be stored as text (matching the existing "Serie" column
# which contains text labels like "01-01-2021") instead of letting Excel
# auto-convert the string into a date serial number.
$ws.Range("A76").NumberFormat = "@"
$ws.Range("A76").Value = "01-04-2021"
$ws.Range("A76").ClearFormats()

$row76Values = @(
    31889,17997,198,152,2215,3275,169,7873,9,54803,
    9044,1840,6263,437,753,19499,5984,1816,9165,10251,
    384,1804,0,302,4,1373,3062,458,2865,72568,
    16925,1762,62,791,24476,1165,302,27083,69533,284,
    10236,223,0,0,12203,12111,300,34176,56444,29685,
    1840,3734,8706,5793,721,1583,4265,116
)

$startCol = 2  # Column B
for ($i = 0; $i -lt $row76Values.Length; $i++) {
    $ws.Cells.Item(76, $startCol + $i).Value = $row76Values[$i]
}
